$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

Set-CellText $ws "D2" "62.105.69"
Set-CellText $ws "E2" "  -0.35%  "
Set-CellText $ws "D3" "2.445.40"
Set-CellText $ws "E3" "  +0.52%  "
Set-CellText $ws "E4" "  -0.11%  "
Set-CellText $ws "D5" "584.45"
Set-CellText $ws "E5" "  +2.49%  "
Set-CellText $ws "D6" "142.78"
Set-CellText $ws "E6" "  -0.48%  "
Set-CellText $ws "E7" "  +0.06%  "
Set-CellText $ws "E8" "  +0.22%  "
Set-CellText $ws "D9" "2.439.37"
Set-CellText $ws "E9" "  +0.43%  "
Set-CellText $ws "E10" "  +1.05%  "
Set-CellText $ws "E11" "  +2.98%  "
Set-CellText $ws "E12" "  +0.22%  "
Set-CellText $ws "E13" "  -2.29%  "
Set-CellText $ws "E14" "  -0.30%  "
Set-CellText $ws "E15" "  +1.44%  "
Set-CellText $ws "D16" "2.887.19"
Set-CellText $ws "E16" "  +0.36%  "
Set-CellText $ws "D17" "62.032.22"
Set-CellText $ws "E17" "  -0.28%  "
Set-CellText $ws "D18" "2.438.54"
Set-CellText $ws "E18" "  +0.37%  "
Set-CellText $ws "D19" "10.76"
Set-CellText $ws "E19" "  -2.49%  "
Set-CellText $ws "D20" "7.14"
Set-CellText $ws "E20" "  +0.19%  "
Set-CellText $ws "D21" "326.28"
Set-CellText $ws "E21" "  +0.53%  "
Set-CellText $ws "D22" "4.10"
Set-CellText $ws "E22" "  -0.69%  "
Set-CellText $ws "E23" "  -0.06%  "
Set-CellText $ws "D24" "1.91"
Set-CellText $ws "E24" "  -5.39%  "
Set-CellText $ws "D25" "65.74"
Set-CellText $ws "E25" "  +1.03%  "
Set-CellText $ws "D26" "9.19"
Set-CellText $ws "E26" "  +1.47%  "
Set-CellText $ws "D27" "594.09"
Set-CellText $ws "E27" "  -4.19%  "
Set-CellText $ws "D28" "0.0₃0972"
Set-CellText $ws "E28" "  +1.14%  "
Set-CellText $ws "D29" "2.566.27"
Set-CellText $ws "E29" "  +0.49%  "
Set-CellText $ws "E30" "  +0.53%  "
Set-CellText $ws "E31" "  -1.80%  "
Set-CellText $ws "D32" "7.98"
Set-CellText $ws "E32" "  -0.57%  "
Set-CellText $ws "E33" "  +1.33%  "
Set-CellText $ws "D34" "0.136"
Set-CellText $ws "E34" "  +0.43%  "
Set-CellText $ws "D35" "4.89"
Set-CellText $ws "E35" "  -2.87%  "
Set-CellText $ws "E36" "  +0.19%  "
Set-CellText $ws "D37" "154.42"
Set-CellText $ws "E37" "  +4.92%  "
Set-CellText $ws "E38" "  -1.69%  "
Set-CellText $ws "E39" "  -0.07%  "
Set-CellText $ws "E40" "  -0.84%  "
Set-CellText $ws "D41" "5.29"
Set-CellText $ws "E41" "  +1.15%  "
Set-CellText $ws "E42" "  +1.87%  "
Set-CellText $ws "E43" "  -1.05%  "
Set-CellText $ws "D45" "2.52"
Set-CellText $ws "E45" "  +1.81%  "
Set-CellText $ws "D46" "141.76"
Set-CellText $ws "E46" "  -2.35%  "
Set-CellText $ws "E47" "  -1.53%  "
Set-CellText $ws "D48" "0.0₆0266"
Set-CellText $ws "E48" "  +20.03%  "
Set-CellText $ws "D49" "0.600"
Set-CellText $ws "E49" "  +1.09%  "
Set-CellText $ws "E50" "  -0.13%  "
Set-CellText $ws "D51" "19.87"
Set-CellText $ws "E51" "  -1.23%  "
